$d = $word.ActiveDocument

$replacements = @(
    @{old = "36×49="; new = "29×81="},
    @{old = "21×30="; new = "90×70="},
    @{old = "27×60="; new = "69×48="},
    @{old = "73×99="; new = "83×87="},
    @{old = "67×77="; new = "23×96="},
    @{old = "50×14="; new = "32×94="},
    @{old = "20×99="; new = "79×79="},
    @{old = "28×79="; new = "92×20="},
    @{old = "62×85="; new = "96×67="},
    @{old = "86×48="; new = "50×85="},
    @{old = "49×81="; new = "46×63="},
    @{old = "45×61="; new = "90×18="},
    @{old = "29×18="; new = "67×23="},
    @{old = "58×95="; new = "35×54="},
    @{old = "69×22="; new = "95×20="},
    @{old = "55×36="; new = "11×34="},
    @{old = "53×81="; new = "95×45="},
    @{old = "68×36="; new = "72×23="},
    @{old = "54×42="; new = "76×77="},
    @{old = "78×70="; new = "33×43="},
    @{old = "49×41="; new = "62×47="},
    @{old = "12×65="; new = "89×81="},
    @{old = "14×36="; new = "69×52="},
    @{old = "64×65="; new = "46×96="},
    @{old = "37×63="; new = "42×11="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
